# Update the BEV vehicles per month workbook with the newly scraped
# registration figures for column Q (2024) and refresh the sheet view
# state left behind by the author's last interactive session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New monthly values for 2024 (column Q), rows 5-8 ---------------------
$ws.Range("Q5").Value = 29668
$ws.Range("Q6").Value = 29708
$ws.Range("Q7").Value = 43412

# Q8 also picks up a thousands-separator number format (new cellXfs entry)
$ws.Range("Q8").Value = 30762
$ws.Range("Q8").NumberFormat = "#,##0"

# --- Sheet view state: scrolled to column C, zoomed to 66%, ---------------
# --- selection moved to Q20 ------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 66
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("Q20").Select()
